$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.1550612449646
$ws.Range("B1").Value = 2.483020067214966
$ws.Range("C1").Value = 4.298581600189209
$ws.Range("D1").Value = 3.499677896499634
$ws.Range("E1").Value = 1.230964303016663
